# Apply "new TPM" data update:
#  - The old row 2 (Target cluster = "ECs") is removed entirely.
#  - The old row 3 (Target cluster = "MuSCs") becomes the new row 2, and its
#    K:T (Receptor-expressing cells .. Edge total expression derived
#    specificity) values are refreshed with newly computed TPM figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the entire row that holds the "ECs" target-cluster record (row 2).
# This shifts the former row 3 ("MuSCs") up to row 2, and since no cell
# references "ECs" anymore, it naturally drops out of the shared-strings
# table when the workbook is saved.
$ws.Rows.Item(2).Delete() | Out-Null

# Update the (now) row 2 values for columns K through T with the new TPM
# derived numbers.
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 0.309635
$ws.Range("N2").Value = 0.928905
$ws.Range("O2").Value = 1
$ws.Range("P2").Value = 1
$ws.Range("Q2").Value = 0.03588215518666667
$ws.Range("R2").Value = 0.32293939668
$ws.Range("S2").Value = 1
$ws.Range("T2").Value = 1
